$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - copy formatting from existing header cell H1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J45
$data = @(
    @(2, 8, 8),
    @(3, 7, 7),
    @(4, 7, 7),
    @(5, 8, 8),
    @(6, 9, 9),
    @(7, 8, 8),
    @(8, 8, 8),
    @(9, 7, 7),
    @(10, 9, 9),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 8, 8),
    @(14, 7, 7),
    @(15, 9, 9),
    @(16, 7, 7),
    @(17, 6, 6),
    @(18, 9, 9),
    @(19, 8, 8),
    @(20, 8, 8),
    @(21, 7, 7),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 8, 8),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 8, 8),
    @(30, 7, 7),
    @(31, 8, 8),
    @(32, 7, 7),
    @(33, 8, 8),
    @(34, 7, 7),
    @(35, 8, 8),
    @(36, 9, 9),
    @(37, 7, 7),
    @(38, 5, 7),
    @(39, 8, 8),
    @(40, 6, 6),
    @(41, 6, 6),
    @(42, 5, 5),
    @(43, 5, 5),
    @(44, 3, 3),
    @(45, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
